$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.037400007247925
$ws.Range("B3").Value = 3.324700117111206
$ws.Range("B4").Value = 0.8774999976158142
$ws.Range("B5").Value = 44.90380096435547
$ws.Range("B6").Value = 29.61039924621582
$ws.Range("B7").Value = 2.078200101852417
$ws.Range("B8").Value = 0.5202999711036682

$ws.Range("A9").Value = "Пастбище"
$ws.Range("B9").Value = 83.35230255126953
$ws.Range("C9").Value = 0.0107
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.2856999933719635
$ws.Range("F9").Value = 0
